$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the B column "BOOL" -> "boolean" for the two boolean-typed vars ---
$ws.Range("B5").Value = "boolean"
$ws.Range("B9").Value = "boolean"

# --- Column C header ---
$ws.Range("C1").Value = "Value"
$ws.Range("C1").Font.Bold = $true

# --- New "Value" column contents for the variable rows ---
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = "E"
$ws.Range("C5").Value = $true
$ws.Range("C6").Value = 10000
$ws.Range("C7").Value = 25
$ws.Range("C8").Value = "A"
$ws.Range("C9").Value = $false

# Left-align the new Value column data cells (C2:C9)
$ws.Range("C2:C9").HorizontalAlignment = -4131

# Match the author's final cursor position
$ws.Range("E18").Select() | Out-Null
